# Modifying view to match new model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old model had A2=123 and E3=1 (dimension A2:E3).
# The new model is a single row A1:C1 = 1, 2, 3.
# Clear the cells that no longer hold data in the new model.
$ws.Range("A2").ClearContents()
$ws.Range("E3").ClearContents()

# Write the new model values.
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3

# Update the selection/active cell to match the new model (C1).
$ws.Range("C1").Select()
